$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053326784846613
$ws.Range("D2").Value = 1.05735288753316
$ws.Range("E2").Value = 1.049784570808108
$ws.Range("F2").Value = 1.065580750290415
$ws.Range("I2").Value = 1.042826644038349
$ws.Range("J2").Value = 1.058344714004051
$ws.Range("K2").Value = 1.060088090727344
$ws.Range("L2").Value = 1.052540665006015
$ws.Range("M2").Value = 1.068293606303547
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.055145191482075
$ws.Range("D3").Value = 1.058813802458063
$ws.Range("E3").Value = 1.051380687837002
$ws.Range("F3").Value = 1.067271853948689
$ws.Range("I3").Value = 1.043339064684925
$ws.Range("J3").Value = 1.059809744085155
$ws.Range("K3").Value = 1.061361181682894
$ws.Range("L3").Value = 1.05394707619614
$ws.Range("M3").Value = 1.069797949821406
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056318172525178
$ws.Range("D4").Value = 1.05975583367735
$ws.Range("E4").Value = 1.052409786857123
$ws.Range("F4").Value = 1.068363033590878
$ws.Range("I4").Value = 1.043667506075708
$ws.Range("J4").Value = 1.060753808169771
$ws.Range("K4").Value = 1.062181170397339
$ws.Range("L4").Value = 1.054852949132864
$ws.Range("M4").Value = 1.070767819387793
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.05681043870908
$ws.Range("D5").Value = 1.060151092597324
$ws.Range("E5").Value = 1.052841551764208
$ws.Range("F5").Value = 1.068821044613688
$ws.Range("I5").Value = 1.043804839804518
$ws.Range("J5").Value = 1.061149772308974
$ws.Range("K5").Value = 1.062525000240248
$ws.Range("L5").Value = 1.055232795306992
$ws.Range("M5").Value = 1.071174719554556
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056893042766477
$ws.Range("D6").Value = 1.06021741349697
$ws.Range("E6").Value = 1.052913996545844
$ws.Range("F6").Value = 1.068897904863724
$ws.Range("I6").Value = 1.043827855344991
$ws.Range("J6").Value = 1.061216202897955
$ws.Range("K6").Value = 1.06258267876948
$ws.Range("L6").Value = 1.055296515960157
$ws.Range("M6").Value = 1.071242991462382
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056324753544351
$ws.Range("D7").Value = 1.059761118156986
$ws.Range("E7").Value = 1.052415559510554
$ws.Range("F7").Value = 1.068369156367585
$ws.Range("I7").Value = 1.043669344047017
$ws.Range("J7").Value = 1.060759102657734
$ws.Range("K7").Value = 1.062185768161253
$ws.Range("L7").Value = 1.054858028496971
$ws.Range("M7").Value = 1.070773259660678
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.05394209086693
$ws.Range("D8").Value = 1.057847297775571
$ws.Range("E8").Value = 1.050324761095038
$ws.Range("F8").Value = 1.066152913502679
$ws.Range("I8").Value = 1.043000470553091
$ws.Range("J8").Value = 1.058840646259902
$ws.Range("K8").Value = 1.060519129628746
$ws.Range("L8").Value = 1.053016840751953
$ws.Range("M8").Value = 1.068802748432776
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.049714739621802
$ws.Range("D9").Value = 1.054449166956319
$ws.Range("E9").Value = 1.046611465790499
$ws.Range("F9").Value = 1.062223324188587
$ws.Range("I9").Value = 1.041797590977813
$ws.Range("J9").Value = 1.055429458626092
$ws.Range("K9").Value = 1.057552712841915
$ws.Range("L9").Value = 1.049739811779812
$ws.Range("M9").Value = 1.065302665906675
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.046875959756139
$ws.Range("D10").Value = 1.052165562565986
$ws.Range("E10").Value = 1.044115371162573
$ws.Range("F10").Value = 1.059586261640353
$ws.Range("I10").Value = 1.040978987272023
$ws.Range("J10").Value = 1.053133775418171
$ws.Range("K10").Value = 1.055554380637495
$ws.Range("L10").Value = 1.047532222864535
$ws.Range("M10").Value = 1.062949653117611
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.045641597991867
$ws.Range("D11").Value = 1.051172225885933
$ws.Range("E11").Value = 1.043029424309392
$ws.Range("F11").Value = 1.058440046503481
$ws.Range("I11").Value = 1.040620478489499
$ws.Range("J11").Value = 1.052134388579281
$ws.Range("K11").Value = 1.054683979983303
$ws.Range("L11").Value = 1.046570662750668
$ws.Range("M11").Value = 1.06192590705435
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.045182304129151
$ws.Range("D12").Value = 1.050802560020184
$ws.Range("E12").Value = 1.042625265688195
$ws.Range("F12").Value = 1.05801361807605
$ws.Range("I12").Value = 1.040486696597368
$ws.Range("J12").Value = 1.051762350543846
$ws.Range("K12").Value = 1.054359890600327
$ws.Range("L12").Value = 1.046212627364843
$ws.Range("M12").Value = 1.061544891313652
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.045280860657811
$ws.Range("D13").Value = 1.050881886394824
$ws.Range("E13").Value = 1.04271199514248
$ws.Range("F13").Value = 1.058105119197538
$ws.Range("I13").Value = 1.040515421255381
$ws.Range("J13").Value = 1.051842191434823
$ws.Range("K13").Value = 1.05442944460467
$ws.Range("L13").Value = 1.046289466806597
$ws.Range("M13").Value = 1.061626654758394
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.045603649009503
$ws.Range("D14").Value = 1.051141683496601
$ws.Range("E14").Value = 1.042996032657095
$ws.Range("F14").Value = 1.058404811641828
$ws.Range("I14").Value = 1.040609432642419
$ws.Range("J14").Value = 1.052103652662797
$ws.Range("K14").Value = 1.054657206737893
$ws.Range("L14").Value = 1.0465410852743
$ws.Range("M14").Value = 1.061894427598679
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.045802423270013
$ws.Range("D15").Value = 1.051301660268451
$ws.Range("E15").Value = 1.043170932305514
$ws.Range("F15").Value = 1.058589372265889
$ws.Range("I15").Value = 1.040667274336175
$ws.Range("J15").Value = 1.05226463816042
$ws.Range("K15").Value = 1.054797434153514
$ws.Range("L15").Value = 1.046695999977918
$ws.Range("M15").Value = 1.062059311233691
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.04695776979082
$ws.Range("D16").Value = 1.052231390338559
$ws.Range("E16").Value = 1.044187332260473
$ws.Range("F16").Value = 1.059662238854433
$ws.Range("I16").Value = 1.041002694386312
$ws.Range("J16").Value = 1.053199987301172
$ws.Range("K16").Value = 1.055612037233869
$ws.Range("L16").Value = 1.047595917597794
$ws.Range("M16").Value = 1.063017491483492
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.04768109282549
$ws.Range("D17").Value = 1.052813363017272
$ws.Range("E17").Value = 1.044823507791792
$ws.Range("F17").Value = 1.060334041436656
$ws.Range("I17").Value = 1.041212005428692
$ws.Range("J17").Value = 1.05378526396752
$ws.Range("K17").Value = 1.056121636891838
$ws.Range("L17").Value = 1.048158883683359
$ws.Range("M17").Value = 1.063617214348593
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.048102499452586
$ws.Range("D18").Value = 1.053152382804493
$ws.Range("E18").Value = 1.045194085449997
$ws.Range("F18").Value = 1.060725474175846
$ws.Range("I18").Value = 1.041333702966754
$ws.Range("J18").Value = 1.054126131915772
$ws.Range("K18").Value = 1.056418385848876
$ws.Range("L18").Value = 1.048486707498538
$ws.Range("M18").Value = 1.063966553335066
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.048246104875314
$ws.Range("D19").Value = 1.053267906518503
$ws.Range("E19").Value = 1.045320359879788
$ws.Range("F19").Value = 1.060858872255759
$ws.Range("I19").Value = 1.041375132751188
$ws.Range("J19").Value = 1.054242272456048
$ws.Range("K19").Value = 1.056519486625203
$ws.Range("L19").Value = 1.048598395111163
$ws.Range("M19").Value = 1.064085589865091
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.047603538523086
$ws.Range("D20").Value = 1.05275096797762
$ws.Range("E20").Value = 1.044755303292519
$ws.Range("F20").Value = 1.06026200678976
$ws.Range("I20").Value = 1.041189588734436
$ws.Range("J20").Value = 1.053722522624628
$ws.Range("K20").Value = 1.056067012667542
$ws.Range("L20").Value = 1.048098539235148
$ws.Range("M20").Value = 1.063552918425883
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.04550861807105
$ws.Range("D21").Value = 1.051065199088407
$ws.Range("E21").Value = 1.042912412648848
$ws.Range("F21").Value = 1.058316578425941
$ws.Range("I21").Value = 1.040581765691919
$ws.Range("J21").Value = 1.052026681671574
$ws.Range("K21").Value = 1.05459015823478
$ws.Range("L21").Value = 1.046467014052367
$ws.Range("M21").Value = 1.061815596047211
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.044186835949537
$ws.Range("D22").Value = 1.050001251478125
$ws.Range("E22").Value = 1.04174913554555
$ws.Range("F22").Value = 1.057089506206535
$ws.Range("I22").Value = 1.040196036733284
$ws.Range("J22").Value = 1.050955678042136
$ws.Range("K22").Value = 1.053657058099633
$ws.Range("L22").Value = 1.045436171283156
$ws.Range("M22").Value = 1.060718918962502
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04488798358006
$ws.Range("D23").Value = 1.050565658769394
$ws.Range("E23").Value = 1.042366251550456
$ws.Range("F23").Value = 1.057740377092088
$ws.Range("I23").Value = 1.040400859581229
$ws.Range("J23").Value = 1.051523894990701
$ws.Range("K23").Value = 1.054152148225783
$ws.Range("L23").Value = 1.045983124456914
$ws.Range("M23").Value = 1.0613007070767
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.047638583504454
$ws.Range("D24").Value = 1.052779162954167
$ws.Range("E24").Value = 1.044786123492072
$ws.Range("F24").Value = 1.060294557438259
$ws.Range("I24").Value = 1.041199719078743
$ws.Range("J24").Value = 1.053750874324706
$ws.Range("K24").Value = 1.056091696522842
$ws.Range("L24").Value = 1.048125807976682
$ws.Range("M24").Value = 1.063581972435238
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.05081114477775
$ws.Range("D25").Value = 1.055330804957632
$ws.Range("E25").Value = 1.047574989343327
$ws.Range("F25").Value = 1.063242198223798
$ws.Range("I25").Value = 1.042111476146532
$ws.Range("J25").Value = 1.05631506035596
$ws.Range("K25").Value = 1.058323193510337
$ws.Range("L25").Value = 1.050590967252725
$ws.Range("M25").Value = 1.06621091043638
